# Update dashboards - 2026-01-29
# Applies the economic-data refresh described by the commit:
#  - shifts "Latest Period" dates forward for several FRED series
#  - updates the Present/Lag1-4 observation columns (Q:U) with newly
#    pulled data
#  - removes the "new data just landed" yellow highlight (style 50)
#    from cells whose date is no longer the most-recently-updated one,
#    restoring them to the plain date style (style 48)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Un-highlight "Latest Period" cells that are no longer fresh.
#    Do this via Copy / PasteSpecial(Formats) from a cell that
#    already carries the plain (non-highlighted) date style, so the
#    workbook reuses the existing style definition instead of
#    minting a new one.
# ---------------------------------------------------------------

$plainDateSrcC = $ws.Range("C24")   # already style 48 (no fill)
$plainDateSrcN = $ws.Range("N3")    # already style 48 (no fill)

$cColumnCells = @("C11","C12","C13","C14","C15","C16","C19","C20","C21","C22","C23")
foreach ($addr in $cColumnCells) {
    $plainDateSrcC.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

$nColumnCells = @("N24","N25","N26","N27","N51")
foreach ($addr in $nColumnCells) {
    $plainDateSrcN.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 2) Refresh "Latest Period" dates (column N) for series whose
#    underlying release date moved forward.
# ---------------------------------------------------------------

$ws.Range("N13").Value = 46041
$ws.Range("N14").Value = 46034
$ws.Range("N29").Value = 46050
$ws.Range("N30").Value = 46050
$ws.Range("N47").Value = 46049
$ws.Range("N48").Value = 46049
$ws.Range("N49").Value = 46049
$ws.Range("N50").Value = 46049
$ws.Range("N52").Value = 46049

# ---------------------------------------------------------------
# 3) Refresh Present / Lag1 / Lag2 / Lag3 / Lag4 observation values.
# ---------------------------------------------------------------

# Row 13 - UI Initial Claims
$ws.Range("Q13").Value = 209000
$ws.Range("R13").Value = 210000
$ws.Range("S13").Value = 199000
$ws.Range("T13").Value = 207000
$ws.Range("U13").Value = 200000

# Row 14 - UI Continuing Claims
$ws.Range("Q14").Value = 1827000
$ws.Range("R14").Value = 1865000
$ws.Range("S14").Value = 1875000
$ws.Range("T14").Value = 1903000
$ws.Range("U14").Value = 1856000

# Row 29 - 5yr, 5yr Forward
$ws.Range("Q29").Value = 2.22
$ws.Range("R29").Value = 2.21
$ws.Range("S29").Value = 2.19
$ws.Range("T29").Value = 2.18
$ws.Range("U29").Value = 2.2

# Row 30 - 10yr TIPS
$ws.Range("Q30").Value = 2.36
$ws.Range("R30").Value = 2.34
$ws.Range("T30").Value = 2.32
$ws.Range("U30").Value = 2.31

# Row 48 - 2y UST
$ws.Range("Q48").Value = 3.53
$ws.Range("R48").Value = 3.56
$ws.Range("S48").Value = 3.6
$ws.Range("T48").Value = 3.61

# Row 49 - 5y UST
$ws.Range("Q49").Value = 3.81
$ws.Range("R49").Value = 3.82
$ws.Range("S49").Value = 3.84
$ws.Range("T49").Value = 3.85
$ws.Range("U49").Value = 3.83

# Row 50 - 10y UST
$ws.Range("Q50").Value = 4.24
$ws.Range("R50").Value = 4.22
$ws.Range("S50").Value = 4.24
$ws.Range("U50").Value = 4.26

# Row 52 - BAA
$ws.Range("Q52").Value = 5.85
$ws.Range("R52").Value = 5.83
$ws.Range("T52").Value = 5.85
$ws.Range("U52").Value = 5.88
